$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row: rename "Username" -> "OpcAuthenticationUsername"
#             and "Password" -> "OpcAuthenticationPassword"
$ws.Range("D1").Value = "OpcAuthenticationUsername"
$ws.Range("E1").Value = "OpcAuthenticationPassword"

# Update the active selection to match the edited cell
$ws.Range("E2").Select()
